$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 86
$lastRow = 106

# Data for columns A (LM id), B (Title), C (Price) - rows 86..106
$dataA = @(
    "91989296",
    "86839655",
    "86839655",
    "85639624",
    "89998902",
    "89998902",
    "89998902",
    "89841815",
    "89841815",
    "89841815",
    "89841822",
    "89841822",
    "89841822",
    "89841822",
    "89841822",
    "89841822",
    "89841822",
    "89841822",
    "89825036",
    "89825036",
    "89380725"
)
$dataB = @(
    "Banheira de Imersão Zen 150x72cm Branco Sensea",
    "Cabo Flexível 2,5mm 100m Azul 750V SIL Fios",
    "Cabo Flexível 2,5mm 100m Azul 750V SIL Fios",
    "Fita Isolante 3M Scotch 33+ Uso Profissional Classe A Preta 19mm x 20m x 0,19mm",
    "Kit Vaso Sanitário com Caixa Acoplada e Assento Branco Duplo Acionamento 3/6L Saída Vertical Gap Roca",
    "Kit Vaso Sanitário com Caixa Acoplada e Assento Branco Duplo Acionamento 3/6L Saída Vertical Gap Roca",
    "Kit Vaso Sanitário com Caixa Acoplada e Assento Branco Duplo Acionamento 3/6L Saída Vertical Gap Roca",
    "Kit Vaso Sanitário com Caixa Acoplada e Assento Branco Duplo Acionamento 3/6L Saída Vertical Lille Sensea",
    "Kit Vaso Sanitário com Caixa Acoplada e Assento Branco Duplo Acionamento 3/6L Saída Vertical Lille Sensea",
    "Kit Vaso Sanitário com Caixa Acoplada e Assento Branco Duplo Acionamento 3/6L Saída Vertical Lille Sensea",
    "Kit Vaso Sanitário com Caixa Acoplada e Assento Branco Duplo Acionamento 3/6L Saída Vertical Lyon Sensea",
    "Kit Vaso Sanitário com Caixa Acoplada e Assento Branco Duplo Acionamento 3/6L Saída Vertical Lyon Sensea",
    "Kit Vaso Sanitário com Caixa Acoplada e Assento Branco Duplo Acionamento 3/6L Saída Vertical Lyon Sensea",
    "Kit Vaso Sanitário com Caixa Acoplada e Assento Branco Duplo Acionamento 3/6L Saída Vertical Lyon Sensea",
    "Kit Vaso Sanitário com Caixa Acoplada e Assento Branco Duplo Acionamento 3/6L Saída Vertical Lyon Sensea",
    "Kit Vaso Sanitário com Caixa Acoplada e Assento Branco Duplo Acionamento 3/6L Saída Vertical Lyon Sensea",
    "Kit Vaso Sanitário com Caixa Acoplada e Assento Branco Duplo Acionamento 3/6L Saída Vertical Lyon Sensea",
    "Kit Vaso Sanitário com Caixa Acoplada e Assento Branco Duplo Acionamento 3/6L Saída Vertical Lyon Sensea",
    "Kit Vaso Sanitário com Caixa Acoplada e Assento Branco Duplo Acionamento 3/6L Saída Vertical Quadra Deca",
    "Kit Vaso Sanitário com Caixa Acoplada e Assento Branco Duplo Acionamento 3/6L Saída Vertical Quadra Deca",
    "Porta Papel Higiênico Metal Dupla Sensea"
)
$dataC = @(
    "8.999.00",
    "159.90",
    "159.90",
    "25.90",
    "1.799.00",
    "1.799.00",
    "1.799.00",
    "804.90",
    "804.90",
    "804.90",
    "899.90",
    "899.90",
    "899.90",
    "899.90",
    "899.90",
    "899.90",
    "899.90",
    "899.90",
    "899.90",
    "899.90",
    "62.90"
)

$rangeA = $ws.Range("A" + $firstRow + ":A" + $lastRow)
$rangeB = $ws.Range("B" + $firstRow + ":B" + $lastRow)
$rangeC = $ws.Range("C" + $firstRow + ":C" + $lastRow)

# Force text storage for numeric-looking values in columns A and C
$rangeA.NumberFormat = "@"
$rangeC.NumberFormat = "@"

for ($i = 0; $i -lt $dataA.Count; $i++) {
    $r = $firstRow + $i
    $ws.Cells.Item($r, 1).Value = $dataA[$i]
    $ws.Cells.Item($r, 2).Value = $dataB[$i]
    $ws.Cells.Item($r, 3).Value = $dataC[$i]
}

# Remove the temporary text formatting so no style is stamped on the new cells
$rangeA.ClearFormats()
$rangeC.ClearFormats()

